$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26
$ws.Range("A26").Value = ""
$ws.Range("B26").Value = "Helle und grosse 4.5 Zimmer Wohnung mit Balkon"
$ws.Range("C26").Value = "Wohnung • 4.5 Zi. • 80 m²"
$ws.Range("D26").Value = "Wohnung • 4.5 Zi. • 80 m²"
$ws.Range("E26").Value = "CHF 1’625."
$ws.Range("F26").Value = "CHF 244"
$ws.Range("H26").Value = "flatfox.ch"

# Row 27
$ws.Range("A27").Value = "MATRO Immobilien AG"
$ws.Range("B27").Value = "4.5 Zimmer-Wohnung mit Balkon"
$ws.Range("C27").Value = "Wohnung • 4.5 Zi. • 96 m²"
$ws.Range("D27").Value = "Wohnung • 4.5 Zi. • 96 m²"
$ws.Range("E27").Value = "CHF 1’590."
$ws.Range("F27").Value = "CHF 199"
$ws.Range("G27").Value = "Bernhardswiesstrasse 29, 9014 St. Gallen"
$ws.Range("H27").Value = "flatfox.ch"

# Row 28
$ws.Range("A28").Value = "MATRO Immobilien AG"
$ws.Range("B28").Value = "3,5 Zimmer-Wohnung mit Sitzplatz"
$ws.Range("C28").Value = "Wohnung • 3.5 Zi. • 75 m²"
$ws.Range("D28").Value = "Wohnung • 3.5 Zi. • 75 m²"
$ws.Range("E28").Value = "CHF 1’430."
$ws.Range("F28").Value = "CHF 229"
$ws.Range("G28").Value = "Bernhardswiesstrasse 31, 9014 St. Gallen"
$ws.Range("H28").Value = "flatfox.ch"

# Row 29
$ws.Range("A29").Value = "René L. Stein"
$ws.Range("B29").Value = "St. Georgen, Parterre, Klusstr. 20"
$ws.Range("C29").Value = "Wohnung • 3.5 Zi. • 71 m²"
$ws.Range("D29").Value = "Wohnung • 3.5 Zi. • 71 m²"
$ws.Range("E29").Value = "CHF 1’735."
$ws.Range("F29").Value = "CHF 293"

# Row 30
$ws.Range("A30").Value = "Sproll & Ramseyer AG"
$ws.Range("B30").Value = "Moderne 1.5-Zimmer-Stadtwohnung im Osten der Stadt"
$ws.Range("C30").Value = "Möblierte Wohnung • 1.5 Zi. • 32 m²"
$ws.Range("D30").Value = "Möblierte Wohnung • 1.5 Zi. • 32 m²"
$ws.Range("E30").Value = "CHF 980."
$ws.Range("F30").Value = "CHF 368"
$ws.Range("G30").Value = "9008 St. Gallen"
$ws.Range("H30").Value = "newhome.ch"

# Row 31
$ws.Range("A31").Value = "Sproll & Ramseyer AG"
$ws.Range("B31").Value = "1.5 Zimmerwohnung direkt an der Uni"
$ws.Range("C31").Value = "Wohnung • 1.5 Zi. • 36 m²"
$ws.Range("D31").Value = "Wohnung • 1.5 Zi. • 36 m²"
$ws.Range("E31").Value = "CHF 980."
$ws.Range("F31").Value = "CHF 327"
$ws.Range("G31").Value = "Höhenweg 15, 9000 St. Gallen"
$ws.Range("H31").Value = "newhome.ch"

# Row 32
$ws.Range("A32").Value = "Wincasa AG"
$ws.Range("B32").Value = "Zentrale Stadtwohnung sucht Sie!"
$ws.Range("C32").Value = "Wohnung • 4 Zi. • 97 m²"
$ws.Range("D32").Value = "Wohnung • 4 Zi. • 97 m²"
$ws.Range("E32").Value = "CHF 2’130."
$ws.Range("F32").Value = "CHF 264"
$ws.Range("G32").Value = "Goliathgasse 37, 9000 St. Gallen"
$ws.Range("H32").Value = "homegate.ch"

# Row 33
$ws.Range("A33").Value = "VTAG Verwaltungs- und Treuhand AG"
$ws.Range("B33").Value = "Bahnhof und Zentrums Nahe Wohnung"
$ws.Range("C33").Value = "Wohnung • 3 Zi."
$ws.Range("D33").Value = "Wohnung • 3 Zi."
$ws.Range("E33").Value = "CHF 1’240."
$ws.Range("F33").Value = "—"
$ws.Range("G33").Value = "Oberstrasse 25, 9000 St. Gallen"
$ws.Range("H33").Value = "immoscout24.ch"

# Row 34
$ws.Range("A34").Value = "Lach Financial Consulting GmbH"
$ws.Range("B34").Value = "Modernisierte 5,5 Zi. Maisonette-Wohnung"
$ws.Range("C34").Value = "Villa • 5.5 Zi. • 133 m²"
$ws.Range("D34").Value = "Villa • 5.5 Zi. • 133 m²"
$ws.Range("E34").Value = "CHF 2’390."
$ws.Range("F34").Value = "CHF 216"
$ws.Range("G34").Value = "9014 St. Gallen"

# Row 35
$ws.Range("A35").Value = "UZON Immobilien AG"
$ws.Range("B35").Value = "Schöne, helle 3.5 Zimmer Wohnung im Grünen"
$ws.Range("C35").Value = "Wohnung • 2.5 Zi."
$ws.Range("D35").Value = "Wohnung • 2.5 Zi."
$ws.Range("E35").Value = "CHF 990."
$ws.Range("F35").Value = "—"
$ws.Range("G35").Value = "Im Hölzli 27, 9000 St. Gallen"
$ws.Range("H35").Value = "immoscout24.ch"

# Row 36
$ws.Range("A36").Value = "UZON Immobilien AG"
$ws.Range("B36").Value = "Gemütliche 1,5 Dachzimmer Wohnung"
$ws.Range("C36").Value = "Dachwohnung • 1.5 Zi."
$ws.Range("D36").Value = "Dachwohnung • 1.5 Zi."
$ws.Range("E36").Value = "CHF 750."
$ws.Range("F36").Value = "—"
$ws.Range("G36").Value = "Sonneggstrasse 5, 9000 St. Gallen"

# Row 37
$ws.Range("A37").Value = "UZON Immobilien AG"
$ws.Range("B37").Value = "Wohnen neben dem Silberturm"
$ws.Range("C37").Value = "Wohnung • 2.5 Zi."
$ws.Range("D37").Value = "Wohnung • 2.5 Zi."
$ws.Range("E37").Value = "CHF 1’300."
$ws.Range("G37").Value = "Sonneggstrasse 5, 9000 St. Gallen"

# Row 38
$ws.Range("A38").Value = "—"
$ws.Range("B38").Value = "Schöne 4-Zimmer-Altbauwohnung in St. Gallen"
$ws.Range("C38").Value = "Wohnung • 4 Zi."
$ws.Range("D38").Value = "Wohnung • 4 Zi."
$ws.Range("E38").Value = "CHF 1’720."
$ws.Range("F38").Value = "—"

# Row 39
$ws.Range("A39").Value = "MATRO Immobilien AG"
$ws.Range("B39").Value = "4,5 Zimmer Wohnung mit hellem Wohnzimmer"
$ws.Range("C39").Value = "Wohnung • 4.5 Zi. • 79 m²"
$ws.Range("D39").Value = "Wohnung • 4.5 Zi. • 79 m²"
$ws.Range("E39").Value = "CHF 1’515."
$ws.Range("F39").Value = "CHF 230"
$ws.Range("G39").Value = "Rorschacher Strasse 233, 9016 St. Gallen"
$ws.Range("H39").Value = "flatfox.ch"

# Row 40
$ws.Range("A40").Value = "Auwiesen Immobilien AG"
$ws.Range("B40").Value = "Wohnen mit Stil und Komfort"
$ws.Range("C40").Value = "Wohnung • 3 Zi. • 92 m²"
$ws.Range("D40").Value = "Wohnung • 3 Zi. • 92 m²"
$ws.Range("E40").Value = "CHF 1’650."
$ws.Range("F40").Value = "CHF 215"
$ws.Range("G40").Value = "9015 St. Gallen"

# Row 41
$ws.Range("A41").Value = "Crowdhouse AG"
$ws.Range("B41").Value = "MODERNE WOHNUNG MIT BALKON SUCHT NACHMIETER"
$ws.Range("C41").Value = "Wohnung • 3 Zi. • 53 m²"
$ws.Range("D41").Value = "Wohnung • 3 Zi. • 53 m²"
$ws.Range("E41").Value = "CHF 1’380."
$ws.Range("F41").Value = "CHF 312"
$ws.Range("G41").Value = "Rickenstrasse 20, 9014 St. Gallen"

# Row 42
$ws.Range("A42").Value = "Crowdhouse AG"
$ws.Range("B42").Value = "ATTRAKTIVE KLEINWOHNUNG IN STADTNÄHE MIT WASCHTURM"
$ws.Range("C42").Value = "Wohnung • 2 Zi. • 43 m²"
$ws.Range("D42").Value = "Wohnung • 2 Zi. • 43 m²"
$ws.Range("E42").Value = "CHF 920."
$ws.Range("F42").Value = "CHF 257"
$ws.Range("G42").Value = "Rickenstrasse 20, 9014 St. Gallen"

# Row 43
$ws.Range("A43").Value = "Dr. Kenig Liegenschaften"
$ws.Range("B43").Value = "MODERNE 2.5 ZIMMERWOHNUNG IN ST. GALLEN"
$ws.Range("C43").Value = "Wohnung • 2.5 Zi."
$ws.Range("D43").Value = "Wohnung • 2.5 Zi."
$ws.Range("E43").Value = "CHF 990."
$ws.Range("F43").Value = "—"

# Row 44
$ws.Range("A44").Value = "—"
$ws.Range("B44").Value = ""
$ws.Range("C44").Value = "Wohnung • 1 Zi. • 22 m²"
$ws.Range("D44").Value = "Wohnung • 1 Zi. • 22 m²"
$ws.Range("E44").Value = "CHF 899."
$ws.Range("F44").Value = "CHF 490"
$ws.Range("G44").Value = "Spisergasse 43, 9000 St. Gallen"

# Row 45
$ws.Range("B45").Value = "Zentrale 1.5 Zimmerwohnung"
$ws.Range("C45").Value = "Wohnung • 1 Zi."
$ws.Range("D45").Value = "Wohnung • 1 Zi."
$ws.Range("E45").Value = "CHF 780."
$ws.Range("G45").Value = "Redingstrasse 8, 9000 St. Gallen"

# Row 46
$ws.Range("A46").Value = "Dafema AG"
$ws.Range("B46").Value = "2-Zimmerwohnung im Krontal"
$ws.Range("C46").Value = "Wohnung • 2 Zi."
$ws.Range("D46").Value = "Wohnung • 2 Zi."
$ws.Range("E46").Value = "CHF 940."
$ws.Range("F46").Value = "—"
$ws.Range("H46").Value = "newhome.ch"

# Row 47
$ws.Range("B47").Value = "Tapetenwechsel gewünscht?"
$ws.Range("C47").Value = "Wohnung • 4.5 Zi. • 100 m²"
$ws.Range("D47").Value = "Wohnung • 4.5 Zi. • 100 m²"
$ws.Range("E47").Value = "CHF 1’750."
$ws.Range("F47").Value = "CHF 210"
$ws.Range("G47").Value = "Oberstrasse 289, 9014 St. Gallen"

# Row 48
$ws.Range("A48").Value = "Brüschweiler Immobilien AG"
$ws.Range("B48").Value = "Grosszügige 6.5 Zimmer Wohnung in St.Gallen, sehr geeignet für ein Büro mit Wohnen"
$ws.Range("C48").Value = "Wohnung • 6.5 Zi."
$ws.Range("D48").Value = "Wohnung • 6.5 Zi."
$ws.Range("E48").Value = "CHF 3’360."
$ws.Range("F48").Value = "—"
$ws.Range("G48").Value = "9016 St. Gallen"
$ws.Range("H48").Value = "icasa.ch"

# Row 49
$ws.Range("A49").Value = "—"
$ws.Range("B49").Value = "2 Zimmer Wohnung Stadtzentrum"
$ws.Range("C49").Value = "Wohnung • 2 Zi."
$ws.Range("D49").Value = "Wohnung • 2 Zi."
$ws.Range("E49").Value = "CHF 855."
$ws.Range("F49").Value = "—"
$ws.Range("H49").Value = "newhome.ch"

# Row 314
$ws.Range("A314").Value = "THOMA Immobilien Treuhand AG"
$ws.Range("B314").Value = "Befristetes Wohnen bis Ende Februar 2025"
$ws.Range("C314").Value = "Wohnung • 4.5 Zi. • 164 m²"
$ws.Range("D314").Value = "Wohnung • 4.5 Zi. • 164 m²"
$ws.Range("E314").Value = "CHF 2’010."
$ws.Range("F314").Value = "CHF 147"
$ws.Range("G314").Value = "Dufourstrasse 114, 9000 St. Gallen"
$ws.Range("H314").Value = "newhome.ch"

# Row 315
$ws.Range("A315").Value = "Furter & Furter AG"
$ws.Range("B315").Value = "Studentenzimmer in Jugendstilvilla zu vermieten (Zimmer 5)"
$ws.Range("C315").Value = "Wohnung • 1 Zi. • 29 m²"
$ws.Range("D315").Value = "Wohnung • 1 Zi. • 29 m²"
$ws.Range("E315").Value = "CHF 880."
$ws.Range("F315").Value = "CHF 364"
$ws.Range("G315").Value = "9000 St. Gallen"
$ws.Range("H315").Value = "newhome.ch"

# Row 316
$ws.Range("A316").Value = "cosyhome ag"
$ws.Range("B316").Value = "St. Gallen-Burggraben"
$ws.Range("C316").Value = "Wohnung • 2 Zi. • 25 m²"
$ws.Range("D316").Value = "Wohnung • 2 Zi. • 25 m²"
$ws.Range("E316").Value = "CHF 785."
$ws.Range("F316").Value = "CHF 377"
$ws.Range("G316").Value = "9000 St. Gallen"

# Row 317
$ws.Range("A317").Value = "Crowdhouse AG"
$ws.Range("B317").Value = "ATTRAKTIVE UND RENOVIERTE 2 ZIMMER WOHNUNG IN STADTNÄHE"
$ws.Range("C317").Value = "Wohnung • 2 Zi. • 43 m²"
$ws.Range("D317").Value = "Wohnung • 2 Zi. • 43 m²"
$ws.Range("E317").Value = "CHF 1’280."
$ws.Range("F317").Value = "CHF 357"
$ws.Range("G317").Value = "Rickenstrasse 20, 9014 St. Gallen"

# Row 318
$ws.Range("A318").Value = "HEV Verwaltungs AG"
$ws.Range("B318").Value = "Sanierte 3.5-Zi. Wohnung"
$ws.Range("C318").Value = "Wohnung • 3.5 Zi. • 60 m²"
$ws.Range("D318").Value = "Wohnung • 3.5 Zi. • 60 m²"
$ws.Range("E318").Value = "CHF 1’480."
$ws.Range("F318").Value = "CHF 296"
$ws.Range("G318").Value = "9000 St. Gallen"

# Row 319
$ws.Range("A319").Value = "HGT Immobilien-Treuhand AG"
$ws.Range("B319").Value = "Zentrale 1.5 Zimmer-Wohnung Nähe Uni/Bahnhof/Stadtzentrum"
$ws.Range("C319").Value = "Wohnung • 1.5 Zi. • 30 m²"
$ws.Range("D319").Value = "Wohnung • 1.5 Zi. • 30 m²"
$ws.Range("E319").Value = "CHF 870."
$ws.Range("F319").Value = "CHF 348"
$ws.Range("G319").Value = "9000 St. Gallen"

# Row 320
$ws.Range("A320").Value = "ImmoLeu AG"
$ws.Range("B320").Value = "Traumhafte 5-Zimmer-Wohnung am Rosenberg mit Panoramablick"
$ws.Range("C320").Value = "Wohnung • 5 Zi. • 150 m²"
$ws.Range("D320").Value = "Wohnung • 5 Zi. • 150 m²"
$ws.Range("E320").Value = "CHF 3’400."
$ws.Range("F320").Value = "CHF 272"
$ws.Range("G320").Value = "9000 St. Gallen"
$ws.Range("H320").Value = "newhome.ch"

# Row 321
$ws.Range("A321").Value = "IBSG AG"
$ws.Range("B321").Value = "MODERNE 1 Zimmer-Wohnung"
$ws.Range("C321").Value = "Wohnung • 1 Zi. • 23 m²"
$ws.Range("D321").Value = "Wohnung • 1 Zi. • 23 m²"
$ws.Range("E321").Value = "CHF 935."
$ws.Range("F321").Value = "CHF 488"
$ws.Range("G321").Value = "Brauerstrasse, 9016 St. Gallen"
$ws.Range("H321").Value = "icasa.ch"

# Row 322
$ws.Range("A322").Value = "Unihome GmbH"
$ws.Range("B322").Value = "3.5 Zimmer Wohnung"
$ws.Range("C322").Value = "Wohnung • 3.5 Zi. • 80 m²"
$ws.Range("D322").Value = "Wohnung • 3.5 Zi. • 80 m²"
$ws.Range("E322").Value = "CHF 1’425."
$ws.Range("F322").Value = "CHF 214"
$ws.Range("G322").Value = "9000 St. Gallen"
$ws.Range("H322").Value = "newhome.ch"

# Row 323
$ws.Range("A323").Value = "Immo10 AG"
$ws.Range("B323").Value = "1.5 Zimmer-Studentenwohnung in Stadtnähe"
$ws.Range("C323").Value = "Wohnung • 1.5 Zi."
$ws.Range("D323").Value = "Wohnung • 1.5 Zi."
$ws.Range("E323").Value = "CHF 831."
$ws.Range("F323").Value = "—"
$ws.Range("H323").Value = "newhome.ch"

# Row 324
$ws.Range("A324").Value = "HN Verwaltungs GmbH"
$ws.Range("B324").Value = "Moderne 4.5 Maisonette/Duplex in St. Gallen"
$ws.Range("C324").Value = "Duplex • 4.5 Zi. • 94 m²"
$ws.Range("D324").Value = "Duplex • 4.5 Zi. • 94 m²"
$ws.Range("E324").Value = "CHF 2’150."
$ws.Range("F324").Value = "CHF 274"
$ws.Range("G324").Value = "9008 St. Gallen"

# Row 325
$ws.Range("A325").Value = "Liegenschaften Treuhand St.Gallen AG"
$ws.Range("B325").Value = "4,5 Zimmer-Wohnung 1.OG zu vermieten"
$ws.Range("C325").Value = "Wohnung • 4.5 Zi. • 95 m²"
$ws.Range("D325").Value = "Wohnung • 4.5 Zi. • 95 m²"
$ws.Range("E325").Value = "CHF 1’510."
$ws.Range("F325").Value = "CHF 191"
$ws.Range("G325").Value = "Feldlistrasse 17, 9000 St. Gallen"

# Row 326
$ws.Range("A326").Value = "—"
$ws.Range("B326").Value = "Appartement à Saint-Gall"
$ws.Range("C326").Value = "Wohnung • 3.5 Zi. • 80 m²"
$ws.Range("D326").Value = "Wohnung • 3.5 Zi. • 80 m²"
$ws.Range("E326").Value = "CHF 1’350."
$ws.Range("F326").Value = "CHF 203"
$ws.Range("G326").Value = "9016 St. Gallen"

# Row 327
$ws.Range("A327").Value = "Martin Ilg"
$ws.Range("B327").Value = "Ab sofort, frisch renoviert - Hübsche 2.0 Zi Wohnung mit Sitzplatz"
$ws.Range("C327").Value = "Studio • 2 Zi. • 35 m²"
$ws.Range("D327").Value = "Studio • 2 Zi. • 35 m²"
$ws.Range("E327").Value = "CHF 850."
$ws.Range("F327").Value = "CHF 291"
$ws.Range("G327").Value = "Teufenerstrasse 129, 9012 St. Gallen"
$ws.Range("H327").Value = "homegate.ch"

# Row 328
$ws.Range("B328").Value = "helle 4 Zimmerwohnung"
$ws.Range("C328").Value = "Wohnung • 4 Zi. • 92 m²"
$ws.Range("D328").Value = "Wohnung • 4 Zi. • 92 m²"
$ws.Range("E328").Value = "CHF 1’300."
$ws.Range("F328").Value = "CHF 170"
$ws.Range("G328").Value = "9014 St. Gallen"

# Row 329
$ws.Range("B329").Value = "TOP Wohnung in der Nähe Kantonsspital"
$ws.Range("C329").Value = "Wohnung • 3.5 Zi. • 62 m²"
$ws.Range("D329").Value = "Wohnung • 3.5 Zi. • 62 m²"
$ws.Range("E329").Value = "CHF 1’345."
$ws.Range("F329").Value = "CHF 260"

# Row 330
$ws.Range("A330").Value = "Regimo St. Gallen AG"
$ws.Range("B330").Value = "Mittlerer Rosenberg - ruhige Lage Nähe UNI"
$ws.Range("C330").Value = "Wohnung • 5 Zi."
$ws.Range("D330").Value = "Wohnung • 5 Zi."
$ws.Range("E330").Value = "CHF 1’790."
$ws.Range("F330").Value = "—"
$ws.Range("G330").Value = "9008 St. Gallen"

# Row 331
$ws.Range("A331").Value = "Regimo St. Gallen AG"
$ws.Range("B331").Value = "Renovierte Jugendstilwohnung in St. Fiden - Nähe Kantonsspital"
$ws.Range("C331").Value = "Wohnung • 4 Zi."
$ws.Range("D331").Value = "Wohnung • 4 Zi."
$ws.Range("E331").Value = "CHF 1’720."
$ws.Range("H331").Value = "homegate.ch"

# Row 332
$ws.Range("A332").Value = "Regimo St. Gallen AG"
$ws.Range("B332").Value = "Wohnung an ruhiger und sonniger Lage"
$ws.Range("C332").Value = "Wohnung • 2 Zi. • 58 m²"
$ws.Range("D332").Value = "Wohnung • 2 Zi. • 58 m²"
$ws.Range("E332").Value = "CHF 1’050."
$ws.Range("F332").Value = "CHF 217"
$ws.Range("G332").Value = "9014 St. Gallen"
$ws.Range("H332").Value = "homegate.ch"

# Row 333
$ws.Range("A333").Value = "aro immo ag"
$ws.Range("B333").Value = "Wohnen an zentraler Lage"
$ws.Range("C333").Value = "Wohnung • 3.5 Zi. • 75 m²"
$ws.Range("D333").Value = "Wohnung • 3.5 Zi. • 75 m²"
$ws.Range("E333").Value = "CHF 1’190."
$ws.Range("F333").Value = "CHF 190"
$ws.Range("G333").Value = "9014 St. Gallen"
$ws.Range("H333").Value = "homegate.ch"

# Row 334
$ws.Range("A334").Value = "Dafema AG"
$ws.Range("B334").Value = "Stadtwohnung im modernen Neubau Wassergasse 53"
$ws.Range("C334").Value = "Wohnung • 4.5 Zi. • 132 m²"
$ws.Range("D334").Value = "Wohnung • 4.5 Zi. • 132 m²"
$ws.Range("E334").Value = "CHF 3’260."
$ws.Range("F334").Value = "CHF 296"
$ws.Range("G334").Value = "Wassergasse 53, 9000 St. Gallen"
$ws.Range("H334").Value = "homegate.ch"

# Row 335
$ws.Range("A335").Value = "—"
$ws.Range("B335").Value = "Frisch renovierte sonnige, ruhige, gepflegte 3.5-Zi-Wohnung"
$ws.Range("C335").Value = "Wohnung • 3.5 Zi. • 64 m²"
$ws.Range("D335").Value = "Wohnung • 3.5 Zi. • 64 m²"
$ws.Range("E335").Value = "CHF 1’290."
$ws.Range("F335").Value = "CHF 242"
$ws.Range("H335").Value = "newhome.ch"

# Row 336
$ws.Range("A336").Value = "—"
$ws.Range("B336").Value = "«Neubau Maisonettewohnung mit Garten an Top Lage!»"
$ws.Range("C336").Value = "Duplex • 2.5 Zi. • 66 m²"
$ws.Range("D336").Value = "Duplex • 2.5 Zi. • 66 m²"
$ws.Range("E336").Value = "CHF 1’795."
$ws.Range("F336").Value = "CHF 326"
$ws.Range("G336").Value = "Seeblickstrasse 9, 9010 St. Gallen"
$ws.Range("H336").Value = "urbanhome.ch"

# Row 337
$ws.Range("B337").Value = "3.5 zimmerwohnung an zentraler Lage in St. Gallen"
$ws.Range("C337").Value = "Wohnung • 3.5 Zi. • 70 m²"
$ws.Range("D337").Value = "Wohnung • 3.5 Zi. • 70 m²"
$ws.Range("E337").Value = "CHF 2’520."
$ws.Range("F337").Value = "CHF 432"
$ws.Range("H337").Value = "flatfox.ch"

Write-Host "Applied 343 cell updates"
